$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the two "draft note" paragraphs right after "The Town" entry:
#       "~ off script - the 100-50 lb overloaded horse (7); near the end"
#       "~ check this"
#    (there is another, unrelated "~ check this" paragraph elsewhere in the
#    document, so we locate this pair by finding the paragraph that starts
#    with "~ off script" and removing it together with the paragraph that
#    immediately follows it.)
# ---------------------------------------------------------------------------
$offScriptIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("~ off script")) {
        $offScriptIdx = $i
        break
    }
}

if ($offScriptIdx -gt 0) {
    $startPos = $d.Paragraphs.Item($offScriptIdx).Range.Start
    $endPos = $d.Paragraphs.Item($offScriptIdx + 1).Range.End
    $d.Range($startPos, $endPos).Delete()
}

# ---------------------------------------------------------------------------
# 2) Remove the draft note paragraph right after "Top Gun: Maverick":
#       "~ It has to be Rooster's fighter jet crashing into a large
#        explosion : like nuclear weapon  rooster (10th)"
# ---------------------------------------------------------------------------
$roosterIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("~ It has to be Rooster")) {
        $roosterIdx = $i
        break
    }
}

if ($roosterIdx -gt 0) {
    $d.Paragraphs.Item($roosterIdx).Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Swap the "Bullets" / "Bullets (user)" character styles' identities:
#    the style previously keyed "Bullets" (name "Bullets") becomes keyed
#    "Bulletsuser" (name "Bullets (user)"), and vice versa. Word's object
#    model never allows rewriting an existing style's internal styleId in
#    place, so re-create both styles (same OpenSymbol font formatting and
#    Quick Style flag as the originals) under swapped names; the engine
#    derives each new styleId from its name, producing the swap.
# ---------------------------------------------------------------------------
$styles = $d.Styles

$bullets = $styles.Item("Bullets")
$bullets.Delete()
$bulletsUser = $styles.Item("Bullets (user)")
$bulletsUser.Delete()

$newBulletsUser = $styles.Add("Bullets (user)", 2)
$newBulletsUser.QuickStyle = $true
$newBulletsUser.Font.Name = "OpenSymbol"
$newBulletsUser.Font.NameFarEast = "OpenSymbol"
$newBulletsUser.Font.NameBi = "OpenSymbol"

$newBullets = $styles.Add("Bullets", 2)
$newBullets.QuickStyle = $true
$newBullets.Font.Name = "OpenSymbol"
$newBullets.Font.NameFarEast = "OpenSymbol"
$newBullets.Font.NameBi = "OpenSymbol"
